# Auto-update draw results: append the 2025-11-06 Pick 4 draw as new row 51.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 51
$target = $ws.Range("A" + $row + ":E" + $row)

# Every existing column in this sheet is stored as plain text (even the
# date-looking and number-looking columns). Force text formatting first so
# "2025-11-06" isn't reinterpreted as a date serial and "251106" isn't
# reinterpreted as a number.
$target.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-11-06"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "251106"
$ws.Cells.Item($row, 4).Value = "9-3-7-3"
$ws.Cells.Item($row, 5).Value = "2025-11-06T21:40:40.325+04:00"

# Restore the default ("Normal") style so the new row matches every other
# data row instead of keeping a one-off "@" number-format style.
$target.Style = "Normal"
